$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply the same formatting as the other header cells (bold, bordered, centered)
# by copying the format from the existing H1 header cell.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for I and J columns (rows 2-12)
$values = @(
    @(4, 7),
    @(5, 8),
    @(4, 6),
    @(8, 8),
    @(1, 5),
    @(1, 5),
    @(1, 4),
    @(5, 7),
    @(5, 7),
    @(1, 2),
    @(1, 1)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
